# Rename "congenital" category to "misc_long_term" across the
# regression_name_list1603-style variable sheets.
#
# Sheets "variables_1612".."variables_1646" keep the "congenital" label in
# cell A4, while sheets "variables_1668".."variables_1688" keep it in cell
# A3 (their row 3 has no extra "health" category above it). Update both
# groups in place.

$wb = $excel.ActiveWorkbook

for ($i = 1612; $i -le 1646; $i++) {
    $name = "variables_$i"
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A4").Value = "misc_long_term"
}

for ($i = 1668; $i -le 1688; $i++) {
    $name = "variables_$i"
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A3").Value = "misc_long_term"
}
